$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the Udoo-quad-specific cells in row 5, keep A5/B5 category labels
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("J5").ClearContents()

# Delete the entire row 6 (Accessories for Procesor / Udoo Accessory kit), shifting rows up
$ws.Rows("6:6").Delete()
